# Apply updated cryptocurrency price/volume data to Sheet1.
# Column D price values that look numeric (e.g. 1.002) are written with a
# leading quote-prefix character so Excel stores them as text (matching the
# original inline-string cells) instead of silently converting them to
# numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''26.251.46'
$ws.Range("E2").Value = '  -0.76%  '

$ws.Range("D3").Value = '''1.676.68'
$ws.Range("E3").Value = '  -1.36%  '

$ws.Range("D4").Value = '''1.002'
$ws.Range("E4").Value = '  -0.90%  '

$ws.Range("D5").Value = '''212.13'
$ws.Range("E5").Value = '  -3.38%  '

$ws.Range("D6").Value = '''0.5285'
$ws.Range("E6").Value = '  -3.54%  '

$ws.Range("D7").Value = '''1.002'
$ws.Range("E7").Value = '  -0.84%  '

$ws.Range("D8").Value = '''0.2659'
$ws.Range("E8").Value = '  -3.36%  '

$ws.Range("D9").Value = '''0.06306'
$ws.Range("E9").Value = '  -2.35%  '

$ws.Range("D10").Value = '''21.42'
$ws.Range("E10").Value = '  -2.98%  '

$ws.Range("D11").Value = '''0.07569'

$ws.Range("D12").Value = '''1.684.61'
$ws.Range("E12").Value = '  -0.89%  '

$ws.Range("D13").Value = '''4.468'
$ws.Range("E13").Value = '  -1.98%  '

$ws.Range("D14").Value = '''0.5630'
$ws.Range("E14").Value = '  -3.67%  '

$ws.Range("D15").Value = '''67.09'
$ws.Range("E15").Value = '  +1.86%  '

$ws.Range("D16").Value = '''0.000008047'
$ws.Range("E16").Value = '  -4.37%  '

$ws.Range("D17").Value = '''25.985.34'
$ws.Range("E17").Value = '  -1.94%  '

$ws.Range("E18").Value = '  -0.76%  '

$ws.Range("D19").Value = '''4.826'
$ws.Range("E19").Value = '  -2.62%  '

$ws.Range("D20").Value = '''188.43'
$ws.Range("E20").Value = '  -1.66%  '

$ws.Range("D21").Value = '''10.42'
$ws.Range("E21").Value = '  -5.24%  '

$ws.Range("D22").Value = '''6.194'
$ws.Range("E22").Value = '  -1.08%  '

$ws.Range("E23").Value = '  -0.78%  '

$ws.Range("D24").Value = '''149.84'
$ws.Range("E24").Value = '  +0.44%  '

$ws.Range("D25").Value = '''0.1254'
$ws.Range("E25").Value = '  -5.23%  '

$ws.Range("D26").Value = '''7.585'
$ws.Range("E26").Value = '  -4.20%  '

$ws.Range("D27").Value = '''16.07'
$ws.Range("E27").Value = '  +1.54%  '

$ws.Range("D28").Value = '''0.06202'
$ws.Range("E28").Value = '  -0.91%  '

$ws.Range("D30").Value = '''1.287'
$ws.Range("E30").Value = '  -3.52%  '

$ws.Range("D31").Value = '''3.505'
$ws.Range("E31").Value = '  -3.10%  '

$ws.Range("D32").Value = '''3.445'
$ws.Range("E32").Value = '  -4.48%  '

$ws.Range("D33").Value = '''1.636'
$ws.Range("E33").Value = '  -3.44%  '

$ws.Range("D34").Value = '''1.003'
$ws.Range("E34").Value = '  -3.72%  '

$ws.Range("D35").Value = '''0.6072'
$ws.Range("E35").Value = '  -1.94%  '

$ws.Range("D36").Value = '''2.408'
$ws.Range("E36").Value = '  -0.25%  '

$ws.Range("D37").Value = '''2.736'
$ws.Range("E37").Value = '  -1.35%  '

$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '''0.01618'
$ws.Range("E38").Value = '  -1.88%  '

$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").Value = '''6.096'
$ws.Range("E39").Value = '  -0.85%  '

$ws.Range("D40").Value = '''1.084.05'
$ws.Range("E40").Value = '  -3.19%  '

$ws.Range("D41").Value = '''0.8708'
$ws.Range("E41").Value = '  -1.10%  '

$ws.Range("D42").Value = '''1.007'
$ws.Range("E42").Value = '  -1.06%  '

$ws.Range("D43").Value = '''99.99'
$ws.Range("E43").Value = '  -1.37%  '

$ws.Range("D44").Value = '''1.824.23'
$ws.Range("E44").Value = '  -1.56%  '

$ws.Range("E45").Value = '  -1.65%  '

$ws.Range("D46").Value = '''56.11'
$ws.Range("E46").Value = '  -2.67%  '

$ws.Range("D47").Value = '''1.002'
$ws.Range("E47").Value = '  -0.54%  '

$ws.Range("D48").Value = '''7.997'
$ws.Range("E48").Value = '  -3.14%  '

$ws.Range("D49").Value = '''0.05229'
$ws.Range("E49").Value = '  -1.15%  '

$ws.Range("D50").Value = '''0.4253'
$ws.Range("E50").Value = '  -1.25%  '

$ws.Range("D51").Value = '''5.985'
$ws.Range("E51").Value = '  -2.92%  '
